$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topics")

# The existing table will be rebuilt once the cell grid has its final shape, so
# detach it first (renaming/resizing ListColumns in place proved unreliable).
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Unlist()

# --- Shift the existing columns (D..G) to their new homes (E, H, I, J) -------
# Use Copy (not Value assignment) so the original cell typing is preserved -
# e.g. the "1.0" in the old Schéma column is stored as text (shared string) and
# must stay text once it lands under the new "Schéma2" header, whereas a plain
# Value assignment of the literal string "1.0" would be reinterpreted as the
# number 1. Processed right-to-left so a source column is always read before a
# later step overwrites it.
$ws.Range("G1:G2").Copy($ws.Range("J1:J2"))
$ws.Range("F1:F2").Copy($ws.Range("I1:I2"))
$ws.Range("E1:E2").Copy($ws.Range("H1:H2"))
$ws.Range("D1:D2").Copy($ws.Range("E1:E2"))

# --- Fill in the brand-new columns --------------------------------------------
$ws.Range("D1").Value = "Type clé"
$ws.Range("D2").Value = "String"
$ws.Range("G1").Value = "Nb partitions"
$ws.Range("I1").Value = "Schéma2"
$ws.Range("F1").Value = "Schéma"
$ws.Range("F2").Value = "monitoring.cpu_v01.avsc"
$ws.Range("G2").Value = 1

# --- Rebuild the table over the new A1:J2 extent ------------------------------
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:J2"), $null, 1)
$newTbl.Name = "Tableau1"
$newTbl.TableStyle = "TableStyleMedium13"

$ws.Cells.Select()
